$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(277)
$p2 = $d.Paragraphs.Item(284)
$r = $d.Range($p1.Range.Start, $p2.Range.End)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7086F5F5" w14:textId="7C6BCCD8" w:rsidR="00AF5C44" w:rsidRDefault="00AF5C44" w:rsidP="00AF5C44"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Plotly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> Plots</w:t></w:r></w:p><w:p w14:paraId="4E9E66BD" w14:textId="391F1B37" w:rsidR="00AF5C44" w:rsidRDefault="00AF5C44" w:rsidP="00AF5C44"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Zipcode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> + Yelp </w:t></w:r><w:r w:rsidRPr="00AF5C44"><w:rPr><w:b/></w:rPr><w:t>[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00AF5C44"><w:t>yelp_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00AF5C44"><w:t>plotly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00AF5C44"><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00AF5C44"><w:t>)</w:t></w:r><w:r w:rsidRPr="00AF5C44"><w:rPr><w:b/></w:rPr><w:t>]</w:t></w:r></w:p><w:p w14:paraId="4E5864E1" w14:textId="59BEDCD4" w:rsidR="00B13121" w:rsidRDefault="00AF5C44" w:rsidP="00B13121"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r w:rsidRPr="00AF5C44"><w:t>Shows all restaurants in the area</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w14:paraId="0A2AA4A2" w14:textId="2319A4B9" w:rsidR="00B13121" w:rsidRPr="004E5197" w:rsidRDefault="00B13121" w:rsidP="00B13121"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r w:rsidRPr="004E5197"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Home Box Plot </w:t></w:r><w:r w:rsidR="004E5197"><w:rPr><w:b/></w:rPr><w:t>[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00822E59"><w:t>homeprices</w:t></w:r><w:r w:rsidR="004E5197"><w:t>_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="004E5197"><w:t>plotly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004E5197"><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="004E5197"><w:t>)</w:t></w:r><w:r w:rsidR="004E5197"><w:rPr><w:b/></w:rPr><w:t>]</w:t></w:r></w:p><w:p w14:paraId="7F6CE448" w14:textId="41CFDD58" w:rsidR="00B13121" w:rsidRPr="00822E59" w:rsidRDefault="00B13121" w:rsidP="00B13121"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r w:rsidRPr="004E5197"><w:t xml:space="preserve">Show the mean median and mode for home values in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="004E5197"><w:t>zipcode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00822E59"><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w14:paraId="5E6B3F1A" w14:textId="09724029" w:rsidR="00822E59" w:rsidRPr="002B2EF9" w:rsidRDefault="00822E59" w:rsidP="00B13121"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:t>Takes 2017 and 2018 mean +- standard deviation for each value</w:t></w:r></w:p><w:p w14:paraId="0273E4E0" w14:textId="649BCC29" w:rsidR="00B13121" w:rsidRPr="002B2EF9" w:rsidRDefault="00B13121" w:rsidP="00B13121"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Rent Box Plot</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>homeprices_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>plotly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>]</w:t></w:r></w:p><w:p w14:paraId="47F43B54" w14:textId="4FEB68F3" w:rsidR="00B13121" w:rsidRPr="002B2EF9" w:rsidRDefault="00B13121" w:rsidP="00B13121"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Show the mean median and mode for rent values in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zipcode</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:t>Takes 2017 and 2018 mean +- standard deviation for each value</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Just change query from home to rent to get proper output </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$r.InsertXML($xml)
Write-Output "Done"
